# Rename the worksheet tab (and the <sheet> entry in workbook.xml)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet 1 - Crash Driver Report"

# Update the title cell (A1): "CrashDriverReport Complete" -> "Crash Driver Report"
$ws.Range("A1").Value = "Crash Driver Report"

# Strip the leading "/ext:CrashDriverInfo" prefix from the Exchange Path column (P)
# for the "Person" class rows (rows 4-11).
$ws.Range("P4").Value  = "nc:Person"
$ws.Range("P5").Value  = "nc:Person/nc:PersonBirthDate"
$ws.Range("P6").Value  = "nc:Person/nc:PersonBirthDate/[nc:DateRepresentation]"
$ws.Range("P7").Value  = "nc:Person/nc:PersonBirthDate/[nc:DateRepresentation]/nc:Date"
$ws.Range("P8").Value  = "nc:Person/nc:PersonName"
$ws.Range("P9").Value  = "nc:Person/nc:PersonName/nc:PersonGivenName"
$ws.Range("P10").Value = "nc:Person/nc:PersonName/nc:PersonMiddleName"
$ws.Range("P11").Value = "nc:Person/nc:PersonName/nc:PersonSurName"
